$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "*PATH_TO_PICTURES*" row (used to be row 2); everything
# below shifts up by one.
$ws.Rows.Item(2).Delete()

# Re-insert a blank row right under the header so the header ("variable" /
# "private name") can span two rows once merged below.
$ws.Rows.Item(2).Insert()

# Add a third column with a human readable explanation for every variable.
$ws.Range("C3").Value = "absolute path to directory where attachments will be saved"
$ws.Range("C4").Value = "sender's email address"
$ws.Range("C5").Value = "sender's host "
$ws.Range("C6").Value = "sender's password for his email address"
$ws.Range("C7").Value = "absolute path to directory where attachments will be archived"
$ws.Range("C8").Value = "absolute path to file with data that will be opened in this notepad"
$ws.Range("C9").Value = "absolute path to empty file 'export.csv'"
$ws.Range("C10").Value = "receiver's email adddress"

# Merge the header cells across the two header rows.
$ws.Range("B1:B2").Merge()
$ws.Range("A1:A2").Merge()

# Give the merged header a box border (split across the two physical rows)
# and center the text both horizontally and vertically.
foreach ($addr in @("A1", "B1")) {
    $c = $ws.Range($addr)
    $c.VerticalAlignment = -4108
    $c.Borders.Item(7).LineStyle = 1
    $c.Borders.Item(8).LineStyle = 1
    $c.Borders.Item(9).LineStyle = -4142
    $c.Borders.Item(10).LineStyle = 1
}
foreach ($addr in @("A2", "B2")) {
    $c = $ws.Range($addr)
    $c.VerticalAlignment = -4108
    $c.Borders.Item(7).LineStyle = 1
    $c.Borders.Item(9).LineStyle = 1
    $c.Borders.Item(10).LineStyle = 1
}

# Restore the view's active selection.
$ws.Range("B12").Select()
